# Reorder supervision by year
#
# The "supervision" table (rows 13-19) listed several thesis supervisions out
# of chronological order. This reorders/fixes the year ("when") column and
# re-sorts those rows so the block is consistent with the rest of the table
# (newest first), while keeping every other row untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supervision")
$ws.Activate()

# --- Row 13 (was: Lina María García Hoyos / 2016-2017) ---
$ws.Range("A13").Value = "Biología"
$ws.Range("B13").Value = "2017 - 2018"
$ws.Range("C13").Value = "Maria Alejandra Abello Mozo  "
$ws.Range("D13").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E13").Value = "\textbf{\textit{Trabajo de grado meritorio}}: \textit{Desarrollo y evaluación de una metodología nueva para manipular las variables del atractivo, dominancia y sexo simultáneamente en fotos de caras humanas con el programa PsychoMorph}"

# --- Row 14 (was: Andrés Castellanos-Chacón / 2017-2018) ---
$ws.Range("A14").Value = "Psicología"
$ws.Range("B14").Value = "2017 - 2018"
$ws.Range("C14").Value = "Cindy Paola Moncada Gómez "
$ws.Range("D14").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E14").Value = "Trabajo de grado: \textit{La voz del sexo casual: ¿las características vocales predicen la disposición al sexo sin compromiso en hombres y mujeres? A}"

# --- Row 15 (was: Angie Liliana Pérez Rodríguez / 2016-2018) ---
$ws.Range("A15").Value = "Psicología"
$ws.Range("B15").Value = "2017 - 2018"
$ws.Range("C15").Value = "Laura Milena Estupiñan Aldana  "
$ws.Range("D15").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E15").Value = "Trabajo de grado: \textit{La voz del sexo casual: ¿las características vocales predicen la disposición al sexo sin compromiso en hombres y mujeres? B}"

# --- Row 16 (was: Lina María Morales Sánchez / 2016-2017) ---
$ws.Range("A16").Value = "Psicología"
$ws.Range("B16").Value = "2016 - 2018"
$ws.Range("C16").Value = "Vanesa Díaz Güiza  "
$ws.Range("D16").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E16").Value = "Trabajo de grado: \textit{La voz del sexo casual: ¿las características vocales predicen la disposición al sexo sin compromiso en hombres y mujeres? C}"

# --- Row 17 (was: Cindy Paola Moncada Gómez / 2016-2017) ---
$ws.Range("A17").Value = "Psicología"
$ws.Range("B17").Value = "2016 - 2018"
$ws.Range("C17").Value = "Lina María García Hoyos  "
$ws.Range("D17").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E17").Value = "Trabajo de grado: \textit{¿Se puede determinar si una persona ha sido infiel a partir de su voz?}"

# --- Row 18 (was: Laura Milena Estupiñan Aldana / 2016-2017) ---
$ws.Range("A18").Value = "Psicología"
$ws.Range("B18").Value = "2016 - 2017"
$ws.Range("C18").Value = "Angie Liliana Pérez Rodríguez  "
$ws.Range("D18").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E18").Value = "Trabajo de grado: \textit{Efectos de los niveles hormonales en la discriminación tonal de mujeres}"

# --- Row 19 (was: Vanesa Díaz Güiza / 2016-2018) ---
$ws.Range("A19").Value = "Psicología"
$ws.Range("B19").Value = "2016 - 2017"
$ws.Range("C19").Value = "Lina María Morales Sánchez "
$ws.Range("D19").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E19").Value = "Trabajo de grado: \textit{Discriminación tonal predice satisfacción con pareja y su inversión parental, en hombres y mujeres}"

# Row 13 now holds the longest entry in the block, so it needs to wrap onto
# more lines; row 14 (now a short "A"-variant entry) goes back to a single
# wrapped line's worth of height, like the rest of the block (rows 15-19
# already are, and stay, single-line height).
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 30

# Column A only ever holds "Psicología"/"Biología"/"Pedagogía Musical" - the
# stale 58-character-wide best-fit from the old data is no longer needed now
# that the column was re-evaluated; shrink it back down.
$ws.Columns.Item(1).ColumnWidth = 16.5

# Reflect where the user was last working in the sheet (scrolled down to the
# reordered block, with C23 selected).
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C23").Select()
